$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 16: Using Your Arcane Powers for Fun and Profit
$ws.Range("H16").Value = 30178.4
$ws.Range("J16").Value = 30178.4
$ws.Range("L16").Value = 30178.4
$ws.Range("N16").Value = -30638.4
# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 3209.7646
$ws.Range("I86").Value = 2738.1428
$ws.Range("J86").Value = 3539.9
$ws.Range("K86").Value = 2738.1428
$ws.Range("L86").Value = 3539.9
$ws.Range("M86").Value = -1615.1428
$ws.Range("N86").Value = -5785.9
# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 3209.7646
$ws.Range("I89").Value = 2738.1428
$ws.Range("J89").Value = 3539.9
$ws.Range("K89").Value = 13690.714
$ws.Range("L89").Value = 17699.5
$ws.Range("M89").Value = -8074.714
$ws.Range("N89").Value = -28931.5
# Row 115: 5-bell Energy
$ws.Range("H115").Value = 308
$ws.Range("I115").Value = 308
$ws.Range("K115").Value = 924
$ws.Range("M115").Value = 643
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2782.614
$ws.Range("I138").Value = 1060.3846
$ws.Range("J138").Value = 3291.4546
$ws.Range("K138").Value = 3181.1538
$ws.Range("L138").Value = 9874.363799999999
$ws.Range("M138").Value = 1958.8462
$ws.Range("N138").Value = -20154.3638

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 4359419
$ws.Range("I32").Value = 4557175
$ws.Range("K32").Value = 4557175
$ws.Range("M32").Value = -4556888
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2077.25
$ws.Range("J61").Value = 3013
$ws.Range("L61").Value = 3013
$ws.Range("N61").Value = -3437
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 437.8889
$ws.Range("I122").Value = 437.125
$ws.Range("K122").Value = 1311.375
$ws.Range("M122").Value = 1138.625
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2077.25
$ws.Range("J136").Value = 3013
$ws.Range("L136").Value = 9039
$ws.Range("N136").Value = -14139

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 7: Thank You for Your Business
$ws.Range("H7").Value = 917.3333
$ws.Range("I7").Value = 700
$ws.Range("J7").Value = 1026
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 1026
$ws.Range("M7").Value = -587
$ws.Range("N7").Value = -1252
# Row 107: The Gold Experience
$ws.Range("H107").Value = 4545.6665
$ws.Range("I107").Value = 3149.6667
$ws.Range("J107").Value = 7337.6665
$ws.Range("K107").Value = 3149.6667
$ws.Range("L107").Value = 7337.6665
$ws.Range("M107").Value = -1229.6667
$ws.Range("N107").Value = -11177.6665

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 12: A Sword in Hand
$ws.Range("H12").Value = 1233.9
$ws.Range("I12").Value = 1039.8334
$ws.Range("J12").Value = 1525
$ws.Range("K12").Value = 1039.8334
$ws.Range("L12").Value = 1525
$ws.Range("M12").Value = -869.8334
$ws.Range("N12").Value = -1865
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2570.4443
$ws.Range("I31").Value = 2627
$ws.Range("J31").Value = 2554.2856
$ws.Range("K31").Value = 2627
$ws.Range("L31").Value = 2554.2856
$ws.Range("M31").Value = -2332
$ws.Range("N31").Value = -3144.2856
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2570.4443
$ws.Range("I34").Value = 2627
$ws.Range("J34").Value = 2554.2856
$ws.Range("K34").Value = 2627
$ws.Range("L34").Value = 2554.2856
$ws.Range("M34").Value = -2425
$ws.Range("N34").Value = -2958.2856
# Row 35: Storm of Swords
$ws.Range("H35").Value = 1612
$ws.Range("I35").Value = 757.5
$ws.Range("K35").Value = 757.5
$ws.Range("M35").Value = -463.5
# Row 107: Built to Last
$ws.Range("H107").Value = 1466.5186
$ws.Range("I107").Value = 675.4666999999999
$ws.Range("K107").Value = 675.4666999999999
$ws.Range("M107").Value = 1244.5333

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On
$ws.Range("H11").Value = 1141
$ws.Range("I11").Value = 1203.7727
$ws.Range("J11").Value = 450.5
$ws.Range("K11").Value = 3611.3181
$ws.Range("L11").Value = 1351.5
$ws.Range("M11").Value = -3471.3181
$ws.Range("N11").Value = -1631.5
# Row 32: Convalescence Precedes Essence
$ws.Range("H32").Value = 946
$ws.Range("J32").Value = 971.5
$ws.Range("L32").Value = 2914.5
$ws.Range("N32").Value = -3480.5
# Row 43: Sole Survivor
$ws.Range("H43").Value = 11234.4
$ws.Range("J43").Value = 12643
$ws.Range("L43").Value = 37929
$ws.Range("N43").Value = -38157
# Row 54: Good Eats in Ishgard
$ws.Range("H54").Value = 11037.375
$ws.Range("J54").Value = 11037.375
$ws.Range("L54").Value = 33112.125
$ws.Range("N54").Value = -34230.125
# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 4583.143
$ws.Range("J80").Value = 4908.25
$ws.Range("L80").Value = 14724.75
$ws.Range("N80").Value = -16596.75
# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 4583.143
$ws.Range("J83").Value = 4908.25
$ws.Range("L83").Value = 44174.25
$ws.Range("N83").Value = -53534.25
# Row 115: Mixology
$ws.Range("H115").Value = 7959.4
$ws.Range("I115").Value = 1657.2
$ws.Range("J115").Value = 11110.5
$ws.Range("K115").Value = 4971.6
$ws.Range("L115").Value = 33331.5
$ws.Range("M115").Value = -3796.6
$ws.Range("N115").Value = -35681.5
# Row 141: Ocean Explosion
$ws.Range("H141").Value = 3465.5
$ws.Range("I141").Value = 2416.9092
$ws.Range("K141").Value = 7250.7276
$ws.Range("M141").Value = -2070.7276

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 10: Let's Talk about Hex
$ws.Range("H10").Value = 6700366.5
$ws.Range("I10").Value = 10050000
$ws.Range("K10").Value = 10050000
$ws.Range("M10").Value = -10049831
# Row 15: The Tusk at Hand
$ws.Range("H15").Value = 10247.25
$ws.Range("J15").Value = 10247.25
$ws.Range("L15").Value = 10247.25
$ws.Range("N15").Value = -10823.25
# Row 27: King for a Day
$ws.Range("H27").Value = 11450
$ws.Range("J27").Value = 15000
$ws.Range("L27").Value = 15000
$ws.Range("N27").Value = -15332
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 10291.25
$ws.Range("I80").Value = 2650
$ws.Range("J80").Value = 11382.857
$ws.Range("K80").Value = 2650
$ws.Range("L80").Value = 11382.857
$ws.Range("M80").Value = -1652
$ws.Range("N80").Value = -13378.857
# Row 81: The Grander Temple
$ws.Range("H81").Value = 10247.25
$ws.Range("J81").Value = 10247.25
$ws.Range("L81").Value = 10247.25
$ws.Range("N81").Value = -12243.25
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 10291.25
$ws.Range("I83").Value = 2650
$ws.Range("J83").Value = 11382.857
$ws.Range("K83").Value = 13250
$ws.Range("L83").Value = 56914.285
$ws.Range("M83").Value = -8258
$ws.Range("N83").Value = -66898.285
# Row 84: Man with a Dragon Earring (L)
$ws.Range("H84").Value = 10247.25
$ws.Range("J84").Value = 10247.25
$ws.Range("L84").Value = 30741.75
$ws.Range("N84").Value = -40725.75
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 10287.777
$ws.Range("J126").Value = 11186.875
$ws.Range("L126").Value = 33560.625
$ws.Range("N126").Value = -38500.625
# Row 132: On Board for Lar
$ws.Range("H132").Value = 4375.25
$ws.Range("I132").Value = 1752.75
$ws.Range("K132").Value = 5258.25
$ws.Range("M132").Value = -2728.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 4: Sole Traders
$ws.Range("H4").Value = 35000
$ws.Range("I4").Value = 35000
$ws.Range("K4").Value = 35000
$ws.Range("M4").Value = -34887
# Row 28: My Sole to Take
$ws.Range("H28").Value = 35000
$ws.Range("I28").Value = 35000
$ws.Range("K28").Value = 35000
$ws.Range("M28").Value = -34768
# Row 37: Quicker than Sand
$ws.Range("H37").Value = 35000
$ws.Range("I37").Value = 35000
$ws.Range("K37").Value = 35000
$ws.Range("M37").Value = -34893
# Row 74: Overall, We Blend In
$ws.Range("H74").Value = 22099
$ws.Range("I74").Value = 22099
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 22099
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -21101
$ws.Range("N74").ClearContents()
# Row 77: Eviction Notice (L)
$ws.Range("H77").Value = 22099
$ws.Range("I77").Value = 22099
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 66297
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -61305
$ws.Range("N77").ClearContents()
# Row 122: Hell on Leather
$ws.Range("H122").Value = 3861.4
$ws.Range("I122").Value = 1342.25
$ws.Range("K122").Value = 4026.75
$ws.Range("M122").Value = -1576.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax
$ws.Range("H107").Value = 1496.9584
$ws.Range("I107").Value = 839.6
$ws.Range("K107").Value = 2518.8
$ws.Range("M107").Value = -598.8000000000002
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1854
$ws.Range("I136").Value = 1612
$ws.Range("K136").Value = 4836
$ws.Range("M136").Value = -2286
